$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (corrections to AgTests / AgPosit columns F & G) ---
$ws.Cells.Item(583, 6).Value2 = 29440
$ws.Cells.Item(583, 7).Value2 = 489

$ws.Cells.Item(600, 6).Value2 = 40558

$ws.Cells.Item(608, 6).Value2 = 46044

$ws.Cells.Item(610, 6).Value2 = 33923

$ws.Cells.Item(614, 6).Value2 = 47690

$ws.Cells.Item(615, 6).Value2 = 36685

$ws.Cells.Item(618, 6).Value2 = 37680
$ws.Cells.Item(618, 7).Value2 = 2654

$ws.Cells.Item(621, 6).Value2 = 55538

$ws.Cells.Item(623, 6).Value2 = 14780
$ws.Cells.Item(623, 7).Value2 = 1550

$ws.Cells.Item(624, 6).Value2 = 50584
$ws.Cells.Item(624, 7).Value2 = 3912

$ws.Cells.Item(625, 6).Value2 = 43236
$ws.Cells.Item(625, 7).Value2 = 3527

$ws.Cells.Item(626, 6).Value2 = 19738
$ws.Cells.Item(626, 7).Value2 = 2068

$ws.Cells.Item(627, 6).Value2 = 33013
$ws.Cells.Item(627, 7).Value2 = 2688

$ws.Cells.Item(628, 6).Value2 = 63506
$ws.Cells.Item(628, 7).Value2 = 4153

$ws.Cells.Item(629, 6).Value2 = 45310
$ws.Cells.Item(629, 7).Value2 = 2867

$ws.Cells.Item(630, 6).Value2 = 45529
$ws.Cells.Item(630, 7).Value2 = 2865

$ws.Cells.Item(631, 6).Value2 = 40769
$ws.Cells.Item(631, 7).Value2 = 2693

# --- Append new daily rows 632-634 ---
$ws.Cells.Item(632, 1).Value2 = 44526
$ws.Cells.Item(632, 2).Value2 = 660386
$ws.Cells.Item(632, 3).Value2 = 26658
$ws.Cells.Item(632, 4).Value2 = 9152
$ws.Cells.Item(632, 5).Value2 = 14228
$ws.Cells.Item(632, 6).Value2 = 39788
$ws.Cells.Item(632, 7).Value2 = 2355

$ws.Cells.Item(633, 1).Value2 = 44527
$ws.Cells.Item(633, 2).Value2 = 667961
$ws.Cells.Item(633, 3).Value2 = 21477
$ws.Cells.Item(633, 4).Value2 = 7575
$ws.Cells.Item(633, 5).Value2 = 14274
$ws.Cells.Item(633, 6).Value2 = 21185
$ws.Cells.Item(633, 7).Value2 = 1732

$ws.Cells.Item(634, 1).Value2 = 44528
$ws.Cells.Item(634, 2).Value2 = 673015
$ws.Cells.Item(634, 3).Value2 = 13582
$ws.Cells.Item(634, 4).Value2 = 5054
$ws.Cells.Item(634, 5).Value2 = 14341
$ws.Cells.Item(634, 6).Value2 = 37257
$ws.Cells.Item(634, 7).Value2 = 1751
